$wb = $excel.ActiveWorkbook

# --- Update the conversion note text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 10.36 = 41710.88 pesos`n✅ 41710.88 pesos = 10.32 = 967.4 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 96.5
$wsTasas.Range("O10").Value = 4025.1
$wsTasas.Range("N12").Value = 4040
$wsTasas.Range("O12").Value = 93.7
